# Add a new "images" column (Q) to the statistics sheet, giving the relative
# path of each team's logo image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("Q1").Value = "images"

# Row 2..9 hold one team per row (A2:A9 = Chelsea, Barcelona, Juventus,
# Marseille, PSG, Bayern, Dortmund, Monaco) - add the matching image path.
$ws.Range("Q2").Value = "../images/chelsea.png"
$ws.Range("Q3").Value = "../images/barca.png"
$ws.Range("Q4").Value = "../images/juventus.png"
$ws.Range("Q5").Value = "../images/marseille.png"
$ws.Range("Q6").Value = "../images/psg.jpeg"
$ws.Range("Q7").Value = "../images/bayern.png"
$ws.Range("Q8").Value = "../images/dortmund.png"
$ws.Range("Q9").Value = "../images/monaco.png"

# Set the width of the new column Q. The target stored OOXML column width is
# ~24.8214285714286 character-units; the COM `ColumnWidth` property is offset
# from that stored width by the standard 0.8333333333333334 padding amount,
# so back that out here to land on the closest achievable stored width.
$ws.Columns.Item(17).ColumnWidth = 24.8214285714286 - 0.8333333333333334

# Update the view: scroll so column I becomes the left-most visible column,
# and move the active selection to Q8 (matches the saved view state).
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("Q8").Select() | Out-Null
